$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4 updates
$ws.Range("G4").Value  = 2.2
$ws.Range("H4").Value  = 2.88
$ws.Range("I4").Value  = 4
$ws.Range("J4").Value  = 1.17
$ws.Range("K4").Value  = 5
$ws.Range("R4").Value  = 2.63
$ws.Range("S4").Value  = 1.44
$ws.Range("U4").Value  = 8.5
$ws.Range("W4").Value  = 21
$ws.Range("X4").Value  = 26
$ws.Range("Y4").Value  = 51
$ws.Range("AD4").Value = 7
$ws.Range("AE4").Value = 17
$ws.Range("AG4").Value = 41

# Row 6 updates
$ws.Range("H6").Value  = 4.25
$ws.Range("I6").Value  = 3.95
$ws.Range("N6").Value  = 1.27
$ws.Range("O6").Value  = 3.1
$ws.Range("R6").Value  = 1.31
$ws.Range("S6").Value  = 2.87
$ws.Range("T6").Value  = 16.5
$ws.Range("U6").Value  = 14
$ws.Range("V6").Value  = 9.25
$ws.Range("W6").Value  = 17.5
$ws.Range("Z6").Value  = 28
$ws.Range("AA6").Value = 10.5
$ws.Range("AB6").Value = 11.5
$ws.Range("AC6").Value = 26
$ws.Range("AE6").Value = 35
$ws.Range("AF6").Value = 15
$ws.Range("AJ6").Value = 100

# Row 14 updates
$ws.Range("G14").Value  = 2.9
$ws.Range("H14").Value  = 3.4
$ws.Range("I14").Value  = 2.18
$ws.Range("K14").Value  = 7.7
$ws.Range("Q14").Value  = 2.77
$ws.Range("R14").Value  = 1.7
$ws.Range("S14").Value  = 2.05
$ws.Range("T14").Value  = 10
$ws.Range("U14").Value  = 15.5
$ws.Range("V14").Value  = 10.5
$ws.Range("W14").Value  = 35
$ws.Range("X14").Value  = 24
$ws.Range("Y14").Value  = 30
$ws.Range("Z14").Value  = 7.7
$ws.Range("AA14").Value = 6.7
$ws.Range("AB14").Value = 13.5
$ws.Range("AD14").Value = 8.5
$ws.Range("AE14").Value = 11.25
$ws.Range("AG14").Value = 21
$ws.Range("AH14").Value = 17
$ws.Range("AI14").Value = 26
